# Insert a new weekly price record as row 6 (shifting existing rows 6-10 down to 7-11)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 6, pushing all following rows down by one
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly record
$ws.Cells.Item(6,1).Value = 1
$ws.Cells.Item(6,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6,3).Value = "Arica y Parinacota"
$ws.Cells.Item(6,4).Value = 44592
$ws.Cells.Item(6,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6,5).Value = 15
$ws.Cells.Item(6,6).Value = 100114007
$ws.Cells.Item(6,7).Value = "Jengibre"
$ws.Cells.Item(6,8).Value = "Sin especificar"
$ws.Cells.Item(6,9).Value = "Primera"
$ws.Cells.Item(6,10).Value = 120
$ws.Cells.Item(6,11).Value = 12000
$ws.Cells.Item(6,12).Value = 13000
$ws.Cells.Item(6,13).Value = 12500
$ws.Cells.Item(6,14).Value = "$/caja 13 kilos"
$ws.Cells.Item(6,15).Value = "Perú"
$ws.Cells.Item(6,16).Value = 962
$ws.Cells.Item(6,17).Value = 13
$ws.Cells.Item(6,18).Value = "Hortaliza"
